# Adding some for loops and file system integration
$wb = $excel.ActiveWorkbook

# --- Remove the stray second data row from Sheet1 (row 7: value 34256354) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows(7).Delete()

# --- Add new Sheet3: a copy of the "Dummy template" header block (rows 1-5) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

# Loop the row/column ranges across to replicate Sheet2's template formatting + values
$copyRanges = @("A1", "B1:F1", "A2:F5")
for ($i = 0; $i -lt $copyRanges.Length; $i++) {
    $rng = $copyRanges[$i]
    $ws2.Range($rng).Copy($ws3.Range($rng))
}

for ($r = 2; $r -le 5; $r++) {
    $ws3.Rows($r).RowHeight = $ws2.Rows($r).RowHeight
}

$ws3.Columns(1).ColumnWidth = $ws2.Columns(1).ColumnWidth
$ws3.Columns("B:F").ColumnWidth = $ws2.Columns(2).ColumnWidth
